$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "27.055.64"
$ws.Range("E2").Value = "  -2.45%  "

Set-TextValue 3 4 "1.865.58"
$ws.Range("E3").Value = "  -2.01%  "

Set-TextValue 4 4 "1.000"
$ws.Range("E4").Value = "  +0.24%  "

Set-TextValue 5 4 "306.37"
$ws.Range("E5").Value = "  -1.85%  "

$ws.Range("E6").Value = "  +0.27%  "

Set-TextValue 7 4 "0.5144"
$ws.Range("E7").Value = "  -1.52%  "

Set-TextValue 8 4 "0.3755"

Set-TextValue 9 4 "0.07166"
$ws.Range("E9").Value = "  -0.90%  "

Set-TextValue 10 4 "0.8889"
$ws.Range("E10").Value = "  -2.00%  "

Set-TextValue 11 4 "20.68"
$ws.Range("E11").Value = "  -2.99%  "

Set-TextValue 12 4 "0.07600"
$ws.Range("E12").Value = "  -0.31%  "

Set-TextValue 13 4 "1.859.85"
$ws.Range("E13").Value = "  -2.63%  "

Set-TextValue 14 4 "5.308"
$ws.Range("E14").Value = "  -2.57%  "

Set-TextValue 15 4 "89.51"
$ws.Range("E15").Value = "  -2.87%  "

Set-TextValue 16 4 "1.001"
$ws.Range("E16").Value = "  +0.34%  "

Set-TextValue 17 4 "0.000008456"
$ws.Range("E17").Value = "  -2.63%  "

Set-TextValue 18 4 "14.08"
$ws.Range("E18").Value = "  -3.06%  "

$ws.Range("E19").Value = "  +0.28%  "

Set-TextValue 20 4 "27.088.18"
$ws.Range("E20").Value = "  -2.48%  "

$ws.Range("E21").Value = "  -2.07%  "

Set-TextValue 22 4 "2.074.44"
$ws.Range("E22").Value = "  -3.33%  "

Set-TextValue 23 4 "10.52"
$ws.Range("E23").Value = "  -2.96%  "

Set-TextValue 24 4 "6.454"
$ws.Range("E24").Value = "  -2.23%  "

Set-TextValue 25 4 "1.843"
$ws.Range("E25").Value = "  -1.36%  "

Set-TextValue 26 4 "147.40"
$ws.Range("E26").Value = "  -3.81%  "

Set-TextValue 27 4 "17.99"
$ws.Range("E27").Value = "  -1.69%  "

Set-TextValue 28 4 "2.117"
$ws.Range("E28").Value = "  -2.19%  "

Set-TextValue 29 4 "112.70"
$ws.Range("E29").Value = "  -1.51%  "

Set-TextValue 30 4 "4.662"
$ws.Range("E30").Value = "  -3.83%  "

Set-TextValue 31 4 "4.703"
$ws.Range("E31").Value = "  -3.51%  "

Set-TextValue 32 4 "0.09105"
$ws.Range("E32").Value = "  +1.25%  "

Set-TextValue 33 4 "0.05130"
$ws.Range("E33").Value = "  -2.72%  "

Set-TextValue 34 4 "3.068"
$ws.Range("E34").Value = "  -3.42%  "

Set-TextValue 35 4 "1.157"
$ws.Range("E35").Value = "  -5.90%  "

Set-TextValue 36 4 "0.7269"
$ws.Range("E36").Value = "  -6.34%  "

Set-TextValue 37 4 "0.02043"
$ws.Range("E37").Value = "  -2.12%  "

Set-TextValue 38 4 "3.044"
$ws.Range("E38").Value = "  -0.75%  "

Set-TextValue 39 4 "2.484"
$ws.Range("E39").Value = "  -5.48%  "

$ws.Range("E40").Value = "  -1.60%  "

Set-TextValue 41 4 "0.5339"
$ws.Range("E41").Value = "  -3.36%  "

Set-TextValue 42 4 "6.553"
$ws.Range("E42").Value = "  -1.83%  "

Set-TextValue 43 4 "116.73"
$ws.Range("E43").Value = "  +1.76%  "

Set-TextValue 44 4 "8.281"
$ws.Range("E44").Value = "  -3.00%  "

Set-TextValue 45 4 "0.1471"
$ws.Range("E45").Value = "  -2.68%  "

Set-TextValue 48 4 "9.978"
$ws.Range("E48").Value = "  -4.90%  "

Set-TextValue 49 4 "1.572"
$ws.Range("E49").Value = "  -2.98%  "

Set-TextValue 50 4 "36.57"

Set-TextValue 51 4 "63.87"
$ws.Range("E51").Value = "  -4.50%  "

# Row 46 and 47: coin order swapped (Decentraland now ranked above PaxDollar)
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue 46 4 "0.4641"
$ws.Range("E46").Value = "  -3.49%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue 47 4 "1.000"
$ws.Range("E47").Value = "  +0.33%  "
